# Update "想去人数" (number of people interested) counts that changed
# between the previous data refresh and the new one (gh-pages data
# regeneration at commit 456a3b4).
#
# Sheet "展览" (sheet1) and sheet "全部类型" (sheet4) both contain rows
# for the same exhibitions, so both need to be updated in lockstep.

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F, per worksheet name.
$updates = @{
    "展览"   = @{ 2 = 4304; 3 = 2437; 7 = 58; 10 = 138; 12 = 1603; 14 = 3350; 15 = 227 }
    "全部类型" = @{ 2 = 4304; 3 = 2437; 8 = 58; 12 = 138; 16 = 1603; 18 = 3350; 19 = 227 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}

$wb.Save()
